$d = $word.ActiveDocument

# The "Socks in the Dark" constraints paragraph currently has its final
# sentence split into two runs with the `_GoBack` bookmark sitting between
# them ("...pick a pair of socks t" | bookmark | "hat match..."). Re-writing
# the whole sentence in one Find/Replace pass re-flows it back into a single
# run and drops that stray bookmark split.
$rng = $d.Content
$rng.Find.Execute( `
    "I want to pick a pair of socks that match and furthermore, pick a pair in all 3 colors. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "I want to pick a pair of socks that match and furthermore, pick a pair in all 3 colors. ", `
    2)

# That paragraph is paragraph 17. Insert a brand-new paragraph right after
# it (inherits the same paragraph/run formatting) to hold the new
# "potential solution" text that's being added to problem 2.
$constraintsPara = $d.Paragraphs(17)
$constraintsPara.Range.InsertParagraphAfter()

$solutionPara = $d.Paragraphs(18)
$solutionPara.Range.Text = "A potential solution would be that I put all of my socks in sections with like colors, so that when I choose, I know where I am reaching. Another solution is to just pair the socks together so that they are not scattered in the drawer. That would be a solution that would solve all of the problems."

# The _GoBack bookmark ends up marking the last edit location, which is now
# the (still empty) paragraph right after the new solution paragraph.
$lastEditPara = $d.Paragraphs(19)
$d.Bookmarks.Add("_GoBack", $lastEditPara.Range)
